$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $oldTimestamp)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Jude Coal Mine, China, M2011, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 7; $row++) {
    $wsData.Range("S$row").Value = $newVersion
}
